# Auto-generated: update static computed price/profit values in Sheets per scheduled data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 183.92308
$ws.Range("J33").Value = 322.25
$ws.Range("L33").Value = 322.25
$ws.Range("N33").Value = -780.25
$ws.Range("H40").Value = 1943.7778
$ws.Range("I40").Value = 1870.7142
$ws.Range("K40").Value = 1870.7142
$ws.Range("M40").Value = -1695.7142
$ws.Range("H62").Value = 5862.5
$ws.Range("I62").Value = 3725
$ws.Range("J62").Value = 8000
$ws.Range("K62").Value = 3725
$ws.Range("L62").Value = 8000
$ws.Range("M62").Value = -3101
$ws.Range("N62").Value = -9248
$ws.Range("H65").Value = 5862.5
$ws.Range("I65").Value = 3725
$ws.Range("J65").Value = 8000
$ws.Range("K65").Value = 18625
$ws.Range("L65").Value = 40000
$ws.Range("M65").Value = -15505
$ws.Range("N65").Value = -46240
$ws.Range("H69").Value = 8500
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 8500
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 25500
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -27248
$ws.Range("H72").Value = 8500
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 8500
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 76500
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -85236
$ws.Range("H80").Value = 6009
$ws.Range("I80").Value = 3913
$ws.Range("K80").Value = 11739
$ws.Range("M80").Value = -10741
$ws.Range("H83").Value = 6009
$ws.Range("I83").Value = 3913
$ws.Range("K83").Value = 35217
$ws.Range("M83").Value = -30225
$ws.Range("H111").Value = 1317.2727
$ws.Range("I111").Value = 570.5714
$ws.Range("K111").Value = 1711.7142
$ws.Range("M111").Value = 1355.2858
$ws.Range("I125").Value = 500001400
$ws.Range("J125").Value = 166668580
$ws.Range("K125").Value = 4500012600
$ws.Range("L125").Value = 1500017220
$ws.Range("M125").Value = -4500010140
$ws.Range("N125").Value = -1500022140
$ws.Range("H135").Value = 772
$ws.Range("I135").Value = 772
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 6948
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -4413
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 1191.5
$ws.Range("I137").Value = 1163.4546
$ws.Range("K137").Value = 3490.3638
$ws.Range("M137").Value = -940.3638000000001
$ws.Range("H138").Value = 2419.9487
$ws.Range("I138").Value = 2780.3333
$ws.Range("J138").Value = 2111.0476
$ws.Range("K138").Value = 8340.999899999999
$ws.Range("L138").Value = 6333.1428
$ws.Range("M138").Value = -3200.999899999999
$ws.Range("N138").Value = -16613.1428
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4690
$ws.Range("I32").Value = 4690
$ws.Range("K32").Value = 4690
$ws.Range("M32").Value = -4403
$ws.Range("H45").Value = 12000
$ws.Range("I45").Value = 12000
$ws.Range("K45").Value = 12000
$ws.Range("M45").Value = -11623
$ws.Range("H134").Value = 64997.668
$ws.Range("J134").Value = 64997.668
$ws.Range("L134").Value = 64997.668
$ws.Range("N134").Value = -75137.66800000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1000
$ws.Range("I20").Value = 1000
$ws.Range("K20").Value = 1000
$ws.Range("M20").Value = -753
$ws.Range("H22").Value = 522.0833
$ws.Range("I22").Value = 427.44446
$ws.Range("K22").Value = 427.44446
$ws.Range("M22").Value = -254.44446
$ws.Range("H94").Value = 1400.7273
$ws.Range("I94").Value = 1367.6666
$ws.Range("K94").Value = 1367.6666
$ws.Range("M94").Value = -916.6666
$ws.Range("H99").Value = 2541.7778
$ws.Range("I99").Value = 2109.625
$ws.Range("J99").Value = 5999
$ws.Range("K99").Value = 2109.625
$ws.Range("L99").Value = 5999
$ws.Range("M99").Value = -611.625
$ws.Range("N99").Value = -8995
$ws.Range("H105").Value = 2868.5715
$ws.Range("I105").Value = 3047.5
$ws.Range("K105").Value = 3047.5
$ws.Range("M105").Value = -1300.5
$ws.Range("H122").Value = 70390
$ws.Range("J122").Value = 70390
$ws.Range("L122").Value = 70390
$ws.Range("N122").Value = -80190
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1770.5
$ws.Range("J31").Value = 2098.75
$ws.Range("L31").Value = 2098.75
$ws.Range("N31").Value = -2688.75
$ws.Range("H34").Value = 1770.5
$ws.Range("J34").Value = 2098.75
$ws.Range("L34").Value = 2098.75
$ws.Range("N34").Value = -2502.75
$ws.Range("H99").Value = 5657.778
$ws.Range("I99").Value = 4989
$ws.Range("J99").Value = 6493.75
$ws.Range("K99").Value = 4989
$ws.Range("L99").Value = 6493.75
$ws.Range("M99").Value = -3491
$ws.Range("N99").Value = -9489.75
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H126").Value = 5657.778
$ws.Range("I126").Value = 4989
$ws.Range("J126").Value = 6493.75
$ws.Range("K126").Value = 14967
$ws.Range("L126").Value = 19481.25
$ws.Range("M126").Value = -12497
$ws.Range("N126").Value = -24421.25
$ws.Range("H134").Value = 2874
$ws.Range("I134").Value = 3108.25
$ws.Range("K134").Value = 9324.75
$ws.Range("M134").Value = -6789.75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 201155.81
$ws.Range("I2").Value = 366804
$ws.Range("J2").Value = 139037.75
$ws.Range("K2").Value = 2200824
$ws.Range("L2").Value = 834226.5
$ws.Range("M2").Value = -2200711
$ws.Range("N2").Value = -834452.5
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()
$ws.Range("H37").Value = 79961.336
$ws.Range("J37").Value = 79961.336
$ws.Range("L37").Value = 239884.008
$ws.Range("N37").Value = -240108.008
$ws.Range("H56").Value = 11897.5
$ws.Range("I56").Value = 11897.5
$ws.Range("K56").Value = 11897.5
$ws.Range("M56").Value = -11367.5
$ws.Range("H117").Value = 13061.125
$ws.Range("I117").Value = 449
$ws.Range("J117").Value = 25673.25
$ws.Range("K117").Value = 1347
$ws.Range("L117").Value = 77019.75
$ws.Range("M117").Value = 2095
$ws.Range("N117").Value = -83903.75
$ws.Range("H124").Value = 2933.3333
$ws.Range("J124").Value = 5900
$ws.Range("L124").Value = 17700
$ws.Range("N124").Value = -27520
$ws.Range("H130").Value = 1828.7142
$ws.Range("I130").Value = 1821
$ws.Range("J130").Value = 1848
$ws.Range("K130").Value = 5463
$ws.Range("L130").Value = 5544
$ws.Range("M130").Value = -443
$ws.Range("N130").Value = -15584
$ws.Range("H132").Value = 2919.2
$ws.Range("I132").Value = 592
$ws.Range("J132").Value = 3501
$ws.Range("K132").Value = 5328
$ws.Range("L132").Value = 31509
$ws.Range("M132").Value = -2798
$ws.Range("N132").Value = -36569
$ws.Range("H134").Value = 6043.091
$ws.Range("I134").Value = 1449
$ws.Range("J134").Value = 18294
$ws.Range("K134").Value = 4347
$ws.Range("L134").Value = 54882
$ws.Range("M134").Value = 723
$ws.Range("N134").Value = -65022
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 250
$ws.Range("I6").Value = 250
$ws.Range("K6").Value = 250
$ws.Range("M6").Value = -137
$ws.Range("H16").Value = 250
$ws.Range("I16").Value = 250
$ws.Range("K16").Value = 250
$ws.Range("M16").Value = 0
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2962.1428
$ws.Range("I46").Value = 2962.1428
$ws.Range("K46").Value = 2962.1428
$ws.Range("M46").Value = -2774.1428
$ws.Range("H61").Value = 2723.1
$ws.Range("I61").Value = 1976.5
$ws.Range("K61").Value = 1976.5
$ws.Range("M61").Value = -1774.5
$ws.Range("H68").Value = 2923.75
$ws.Range("J68").Value = 2923.75
$ws.Range("L68").Value = 2923.75
$ws.Range("N68").Value = -4421.75
$ws.Range("H71").Value = 2923.75
$ws.Range("J71").Value = 2923.75
$ws.Range("L71").Value = 14618.75
$ws.Range("N71").Value = -22106.75
$ws.Range("H93").Value = 1963.5714
$ws.Range("I93").Value = 1991
$ws.Range("J93").Value = 1927
$ws.Range("K93").Value = 1991
$ws.Range("L93").Value = 1927
$ws.Range("M93").Value = -743
$ws.Range("N93").Value = -4423
$ws.Range("H113").Value = 2723.1
$ws.Range("I113").Value = 1976.5
$ws.Range("K113").Value = 1976.5
$ws.Range("M113").Value = 193.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 2257500
$ws.Range("J4").Value = 2515000
$ws.Range("L4").Value = 2515000
$ws.Range("N4").Value = -2515226
$ws.Range("H75").Value = 90118
$ws.Range("I75").Value = 90118
$ws.Range("K75").Value = 90118
$ws.Range("M75").Value = -89182
$ws.Range("H78").Value = 90118
$ws.Range("I78").Value = 90118
$ws.Range("K78").Value = 270354
$ws.Range("M78").Value = -265674
$ws.Range("H96").Value = 1838.8
$ws.Range("I96").Value = 1838.8
$ws.Range("K96").Value = 1838.8
$ws.Range("M96").Value = -465.8
$ws.Range("H122").Value = 1424.5
$ws.Range("I122").Value = 1316
$ws.Range("J122").Value = 1750
$ws.Range("K122").Value = 3948
$ws.Range("L122").Value = 5250
$ws.Range("M122").Value = -1498
$ws.Range("N122").Value = -10150
$ws.Range("H126").Value = 6314.154
$ws.Range("I126").Value = 5211.2856
$ws.Range("K126").Value = 15633.8568
$ws.Range("M126").Value = -13163.8568
